$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.536.21"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = "'2.112.00"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.97%  '

$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.59%  '

$ws.Range("D5").Value = "'335.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.90%  '

$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("E7").Value = '  +0.77%  '

$ws.Range("D8").Value = "'0.4557"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.49%  '

$ws.Range("D9").Value = "'54.95"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.88%  '

$ws.Range("D10").Value = "'0.09114"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.18%  '

$ws.Range("D11").Value = "'1.171"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.73%  '

$ws.Range("D12").Value = "'24.56"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").Value = "'2.117.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.52%  '

$ws.Range("D14").Value = "'6.844"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.37%  '

$ws.Range("D15").Value = "'8.115"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +5.83%  '

$ws.Range("E16").Value = '  +5.06%  '

$ws.Range("D17").Value = "'97.01"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.38%  '

$ws.Range("D18").Value = "'1.010"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.58%  '

$ws.Range("D19").Value = "'0.06692"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.61%  '

$ws.Range("E20").Value = '  +0.89%  '

$ws.Range("D21").Value = "'1.008"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.62%  '

$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").Value = "'30.612.18"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.43%  '

$ws.Range("E24").Value = '  +4.47%  '

$ws.Range("D25").Value = "'2.352"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("D26").Value = "'2.359.51"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.18%  '

$ws.Range("D27").Value = "'22.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.49%  '

$ws.Range("D28").Value = "'163.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.56%  '

$ws.Range("D29").Value = "'2.529"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.25%  '

$ws.Range("D30").Value = "'133.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.83%  '

$ws.Range("E31").Value = '  +2.26%  '

$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").Value = "'1.642"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").Value = "'6.364"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.32%  '

$ws.Range("D35").Value = "'3.954"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.22%  '

$ws.Range("D36").Value = "'10.53"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.53%  '

$ws.Range("D37").Value = "'5.901"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +8.39%  '

$ws.Range("D38").Value = "'0.02617"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.56%  '

$ws.Range("D39").Value = "'0.06808"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("D40").Value = "'0.2323"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.25%  '

$ws.Range("D41").Value = "'12.57"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("D42").Value = "'0.6865"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").Value = "'1.258"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.44%  '

$ws.Range("D44").Value = "'14.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.25%  '

$ws.Range("D45").Value = "'0.6441"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.78%  '

$ws.Range("D46").Value = "'2.307"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.31%  '

$ws.Range("D47").Value = "'3.687"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.68%  '

$ws.Range("D48").Value = "'0.00000000358"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +20.04%  '

$ws.Range("E49").Value = '  +0.79%  '

$ws.Range("D50").Value = "'83.18"
$ws.Range("D50").ClearFormats()

$ws.Range("D51").Value = "'0.3344"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +12.28%  '
